# daily auto push: 2025-10-04 13:28 UTC
# Append the day's new record as row 61 (A1:D60 -> A1:D61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
if ($row -lt 61) { $row = 61 }

# Column A holds plain text dates (e.g. "2025/10/04"), not real Excel dates.
# Force text formatting before assigning the value so COM doesn't coerce the
# string into a date serial number, then restore the default "Normal" style
# (matching every other data row, which carries no explicit style) once the
# text value is safely stored.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/04"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = "土"
$ws.Cells.Item($row, 3).Value = 20
$ws.Cells.Item($row, 4).Value = 201
